$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "time_taken" in F1, matching the style of the other header cells (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Per-row "time_taken" timestamp values (plain text, no special style) for F2:F39
$timestamps = @(
  "2021-10-05 10:50:18.031079",
  "2021-10-05 10:50:18.031090",
  "2021-10-05 10:50:18.031094",
  "2021-10-05 10:50:18.031096",
  "2021-10-05 10:50:18.031099",
  "2021-10-05 10:50:18.031102",
  "2021-10-05 10:50:18.031104",
  "2021-10-05 10:50:18.031107",
  "2021-10-05 10:50:18.031110",
  "2021-10-05 10:50:18.031112",
  "2021-10-05 10:50:18.031115",
  "2021-10-05 10:50:18.031118",
  "2021-10-05 10:50:18.031120",
  "2021-10-05 10:50:18.031123",
  "2021-10-05 10:50:18.031125",
  "2021-10-05 10:50:18.031128",
  "2021-10-05 10:50:18.031130",
  "2021-10-05 10:50:18.031133",
  "2021-10-05 10:50:18.031135",
  "2021-10-05 10:50:18.031138",
  "2021-10-05 10:50:18.031140",
  "2021-10-05 10:50:18.031143",
  "2021-10-05 10:50:18.031146",
  "2021-10-05 10:50:18.031148",
  "2021-10-05 10:50:18.031152",
  "2021-10-05 10:50:18.031154",
  "2021-10-05 10:50:18.031157",
  "2021-10-05 10:50:18.031159",
  "2021-10-05 10:50:18.031162",
  "2021-10-05 10:50:18.031164",
  "2021-10-05 10:50:18.031167",
  "2021-10-05 10:50:18.031169",
  "2021-10-05 10:50:18.031172",
  "2021-10-05 10:50:18.031175",
  "2021-10-05 10:50:18.031177",
  "2021-10-05 10:50:18.031180",
  "2021-10-05 10:50:18.031182",
  "2021-10-05 10:50:18.031185"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
